$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AC, AD, AE
$ws.Range("AC1").Value = "Wins"
$ws.Range("AD1").Value = "Losses"
$ws.Range("AE1").Value = "Ties"

# Match the existing header formatting (bold, centered, bordered) used by A1:AB1
$ws.Range("A1").Copy()
$ws.Range("AC1:AE1").PasteSpecial(-4122)

# Data rows 2-38: season record values for every player row
for ($r = 2; $r -le 38; $r++) {
    $ws.Range("AC$r").Value = 85
    $ws.Range("AD$r").Value = 77
    $ws.Range("AE$r").Value = 0
}

Write-Host "Season record columns added"
